# Apply 2020-08-28 data update to Fonds de solidarite volet 1 dataset.
# Updates columns C (nombre_aides) and D (montant_total) for the rows
# that received refreshed figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 3).Value = 38533
$ws.Cells.Item(2, 4).Value = 55715930
$ws.Cells.Item(3, 3).Value = 92478
$ws.Cells.Item(3, 4).Value = 135554551
$ws.Cells.Item(4, 3).Value = 31633
$ws.Cells.Item(4, 4).Value = 46844320
$ws.Cells.Item(5, 3).Value = 8867
$ws.Cells.Item(5, 4).Value = 13177897
$ws.Cells.Item(6, 3).Value = 2049
$ws.Cells.Item(6, 4).Value = 3045471
$ws.Cells.Item(7, 3).Value = 167
$ws.Cells.Item(7, 4).Value = 245593
$ws.Cells.Item(12, 3).Value = 41994
$ws.Cells.Item(12, 4).Value = 56952605
$ws.Cells.Item(13, 3).Value = 9838
$ws.Cells.Item(13, 4).Value = 14230293
$ws.Cells.Item(14, 3).Value = 26339
$ws.Cells.Item(14, 4).Value = 38618877
$ws.Cells.Item(15, 3).Value = 8419
$ws.Cells.Item(15, 4).Value = 12493978
$ws.Cells.Item(16, 3).Value = 2191
$ws.Cells.Item(16, 4).Value = 3256430
$ws.Cells.Item(17, 3).Value = 426
$ws.Cells.Item(17, 4).Value = 628123
$ws.Cells.Item(20, 3).Value = 10368
$ws.Cells.Item(20, 4).Value = 13718148
$ws.Cells.Item(21, 3).Value = 13613
$ws.Cells.Item(21, 4).Value = 19646266
$ws.Cells.Item(22, 3).Value = 32091
$ws.Cells.Item(22, 4).Value = 47091784
$ws.Cells.Item(23, 3).Value = 10350
$ws.Cells.Item(23, 4).Value = 15385310
$ws.Cells.Item(24, 3).Value = 2675
$ws.Cells.Item(24, 4).Value = 3977271
$ws.Cells.Item(27, 3).Value = 11868
$ws.Cells.Item(27, 4).Value = 15842822
$ws.Cells.Item(28, 3).Value = 7822
$ws.Cells.Item(28, 4).Value = 11323007
$ws.Cells.Item(29, 3).Value = 22877
$ws.Cells.Item(29, 4).Value = 33578815
$ws.Cells.Item(30, 3).Value = 7910
$ws.Cells.Item(30, 4).Value = 11765091
$ws.Cells.Item(31, 3).Value = 2000
$ws.Cells.Item(31, 4).Value = 2984419
$ws.Cells.Item(32, 3).Value = 374
$ws.Cells.Item(32, 4).Value = 558415
$ws.Cells.Item(33, 3).Value = 30
$ws.Cells.Item(33, 4).Value = 44893
$ws.Cells.Item(34, 3).Value = 8439
$ws.Cells.Item(34, 4).Value = 11147120
$ws.Cells.Item(35, 3).Value = 3333
$ws.Cells.Item(35, 4).Value = 4814652
$ws.Cells.Item(36, 3).Value = 7984
$ws.Cells.Item(36, 4).Value = 11659821
$ws.Cells.Item(37, 3).Value = 3222
$ws.Cells.Item(37, 4).Value = 4775961
$ws.Cells.Item(38, 3).Value = 836
$ws.Cells.Item(38, 4).Value = 1245223
$ws.Cells.Item(41, 3).Value = 2519
$ws.Cells.Item(41, 4).Value = 3403422
$ws.Cells.Item(42, 3).Value = 17596
$ws.Cells.Item(42, 4).Value = 25444606
$ws.Cells.Item(43, 3).Value = 51954
$ws.Cells.Item(43, 4).Value = 76152670
$ws.Cells.Item(44, 3).Value = 19242
$ws.Cells.Item(44, 4).Value = 28577379
$ws.Cells.Item(45, 3).Value = 5701
$ws.Cells.Item(45, 4).Value = 8486760
$ws.Cells.Item(46, 3).Value = 1238
$ws.Cells.Item(46, 4).Value = 1847545
$ws.Cells.Item(50, 3).Value = 17028
$ws.Cells.Item(50, 4).Value = 22631338
$ws.Cells.Item(51, 3).Value = 2106
$ws.Cells.Item(51, 4).Value = 3054878
$ws.Cells.Item(52, 3).Value = 7143
$ws.Cells.Item(52, 4).Value = 10496642
$ws.Cells.Item(57, 3).Value = 7247
$ws.Cells.Item(57, 4).Value = 9962662
$ws.Cells.Item(58, 3).Value = 1118
$ws.Cells.Item(58, 4).Value = 1848544
$ws.Cells.Item(59, 3).Value = 2754
$ws.Cells.Item(59, 4).Value = 4550588
$ws.Cells.Item(60, 3).Value = 1079
$ws.Cells.Item(60, 4).Value = 1779838
$ws.Cells.Item(61, 3).Value = 372
$ws.Cells.Item(61, 4).Value = 618883
$ws.Cells.Item(64, 3).Value = 1638
$ws.Cells.Item(64, 4).Value = 2520879
$ws.Cells.Item(65, 3).Value = 15684
$ws.Cells.Item(65, 4).Value = 22652021
$ws.Cells.Item(66, 3).Value = 45431
$ws.Cells.Item(66, 4).Value = 66469834
$ws.Cells.Item(67, 3).Value = 15914
$ws.Cells.Item(67, 4).Value = 23644820
$ws.Cells.Item(68, 3).Value = 4630
$ws.Cells.Item(68, 4).Value = 6896051
$ws.Cells.Item(69, 3).Value = 957
$ws.Cells.Item(69, 4).Value = 1423668
$ws.Cells.Item(73, 3).Value = 15316
$ws.Cells.Item(73, 4).Value = 20177052
$ws.Cells.Item(74, 3).Value = 53599
$ws.Cells.Item(74, 4).Value = 78003489
$ws.Cells.Item(75, 3).Value = 150761
$ws.Cells.Item(75, 4).Value = 222101713
$ws.Cells.Item(76, 3).Value = 65149
$ws.Cells.Item(76, 4).Value = 97079663
$ws.Cells.Item(77, 3).Value = 20858
$ws.Cells.Item(77, 4).Value = 31167322
$ws.Cells.Item(78, 3).Value = 4984
$ws.Cells.Item(78, 4).Value = 7444403
$ws.Cells.Item(79, 3).Value = 276
$ws.Cells.Item(79, 4).Value = 409170
$ws.Cells.Item(80, 3).Value = 22
$ws.Cells.Item(80, 4).Value = 31905
$ws.Cells.Item(85, 3).Value = 52866
$ws.Cells.Item(85, 4).Value = 71854108
$ws.Cells.Item(86, 3).Value = 4732
$ws.Cells.Item(86, 4).Value = 6858111
$ws.Cells.Item(87, 3).Value = 11807
$ws.Cells.Item(87, 4).Value = 17343602
$ws.Cells.Item(88, 3).Value = 3940
$ws.Cells.Item(88, 4).Value = 5872583
$ws.Cells.Item(90, 3).Value = 291
$ws.Cells.Item(90, 4).Value = 434012
$ws.Cells.Item(93, 3).Value = 5527
$ws.Cells.Item(93, 4).Value = 7428878
$ws.Cells.Item(94, 3).Value = 1640
$ws.Cells.Item(94, 4).Value = 2362862
$ws.Cells.Item(95, 3).Value = 5298
$ws.Cells.Item(95, 4).Value = 7804378
$ws.Cells.Item(96, 3).Value = 1970
$ws.Cells.Item(96, 4).Value = 2933326
$ws.Cells.Item(98, 3).Value = 196
$ws.Cells.Item(98, 4).Value = 294613
$ws.Cells.Item(101, 3).Value = 3654
$ws.Cells.Item(101, 4).Value = 4834983
$ws.Cells.Item(102, 3).Value = 697
$ws.Cells.Item(102, 4).Value = 1141825
$ws.Cells.Item(103, 3).Value = 419
$ws.Cells.Item(103, 4).Value = 703097
$ws.Cells.Item(104, 3).Value = 154
$ws.Cells.Item(104, 4).Value = 255520
$ws.Cells.Item(107, 3).Value = 11008
$ws.Cells.Item(107, 4).Value = 15969606
$ws.Cells.Item(108, 3).Value = 29625
$ws.Cells.Item(108, 4).Value = 43512550
$ws.Cells.Item(109, 3).Value = 9916
$ws.Cells.Item(109, 4).Value = 14743205
$ws.Cells.Item(110, 3).Value = 2731
$ws.Cells.Item(110, 4).Value = 4071580
$ws.Cells.Item(111, 3).Value = 502
$ws.Cells.Item(111, 4).Value = 748046
$ws.Cells.Item(114, 3).Value = 9948
$ws.Cells.Item(114, 4).Value = 13135684
$ws.Cells.Item(115, 3).Value = 31046
$ws.Cells.Item(115, 4).Value = 44765595
$ws.Cells.Item(116, 3).Value = 67128
$ws.Cells.Item(116, 4).Value = 98226534
$ws.Cells.Item(117, 3).Value = 21645
$ws.Cells.Item(117, 4).Value = 32168038
$ws.Cells.Item(118, 3).Value = 6137
$ws.Cells.Item(118, 4).Value = 9143521
$ws.Cells.Item(120, 3).Value = 84
$ws.Cells.Item(120, 4).Value = 123420
$ws.Cells.Item(124, 3).Value = 26226
$ws.Cells.Item(124, 4).Value = 35004965
$ws.Cells.Item(125, 3).Value = 36767
$ws.Cells.Item(125, 4).Value = 53055956
$ws.Cells.Item(126, 3).Value = 78100
$ws.Cells.Item(126, 4).Value = 114192411
$ws.Cells.Item(127, 3).Value = 24183
$ws.Cells.Item(127, 4).Value = 35892784
$ws.Cells.Item(128, 3).Value = 6497
$ws.Cells.Item(128, 4).Value = 9655358
$ws.Cells.Item(129, 3).Value = 1270
$ws.Cells.Item(129, 4).Value = 1888311
$ws.Cells.Item(131, 3).Value = 18
$ws.Cells.Item(131, 4).Value = 27000
$ws.Cells.Item(133, 3).Value = 32327
$ws.Cells.Item(133, 4).Value = 42905844
$ws.Cells.Item(134, 3).Value = 13533
$ws.Cells.Item(134, 4).Value = 19590446
$ws.Cells.Item(135, 3).Value = 32811
$ws.Cells.Item(135, 4).Value = 48186690
$ws.Cells.Item(136, 3).Value = 11630
$ws.Cells.Item(136, 4).Value = 17279087
$ws.Cells.Item(137, 3).Value = 3006
$ws.Cells.Item(137, 4).Value = 4480241
$ws.Cells.Item(138, 3).Value = 512
$ws.Cells.Item(138, 4).Value = 761990
$ws.Cells.Item(141, 3).Value = 10967
$ws.Cells.Item(141, 4).Value = 14617685
$ws.Cells.Item(142, 3).Value = 35875
$ws.Cells.Item(142, 4).Value = 51815194
$ws.Cells.Item(143, 3).Value = 82788
$ws.Cells.Item(143, 4).Value = 121285929
$ws.Cells.Item(144, 3).Value = 24728
$ws.Cells.Item(144, 4).Value = 36736822
$ws.Cells.Item(145, 3).Value = 6491
$ws.Cells.Item(145, 4).Value = 9686067
$ws.Cells.Item(146, 3).Value = 1470
$ws.Cells.Item(146, 4).Value = 2187230
$ws.Cells.Item(149, 3).Value = 29704
$ws.Cells.Item(149, 4).Value = 40041316
